$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 140 (shifts existing rows 140-158 down to 141-159)
$ws.Rows.Item(140).Insert()

# Populate the new row 140 with the new price entry
$ws.Range("A140").Value = 4
$ws.Range("B140").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C140").Value = "Los Lagos"
$ws.Range("D140").Value = 45077
$ws.Range("E140").Value = 10
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100104
$ws.Range("H140").Value = "Frutos de pepita"
$ws.Range("I140").Value = 100104003
$ws.Range("J140").Value = "Membrillo"
$ws.Range("K140").Value = "Champion"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 80
$ws.Range("N140").Value = 13000
$ws.Range("O140").Value = 14000
$ws.Range("P140").Value = 13500
$ws.Range("Q140").Value = '$/caja 18 kilos empedrada'
$ws.Range("R140").Value = "Región de O'Higgins"
$ws.Range("S140").Value = 750
$ws.Range("T140").Value = 18
